$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the SQL queries in cells C2, B2, B3, B4, B5, B6, B7 ---
# These cells contain SQL text whose JOIN conditions used generic ".id"
# columns; the commit renames them to the qualified "*_id" columns.
$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cells) {
    $range = $ws.Range($addr)
    $text = $range.Value()

    $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    $range.Value = $text
}

# --- Cosmetic sheet tweaks from the same commit ---
# Selection moved from D2 to B2
$ws.Range("B2").Select()

# Column C widened (and no longer "best fit")
$ws.Columns.Item(3).ColumnWidth = 70.5
